# Applies the cryptos-list price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '29.136.61'
Set-TextCell 'E2' '  +1.31%  '
Set-TextCell 'E3' '  +1.66%  '
Set-TextCell 'E4' '  -0.14%  '
Set-TextCell 'D5' '327.05'
Set-TextCell 'E5' '  +0.80%  '
Set-TextCell 'E6' '  -0.07%  '
Set-TextCell 'D7' '0.4612'
Set-TextCell 'E7' '  -0.26%  '
Set-TextCell 'D8' '0.3928'
Set-TextCell 'E8' '  +1.56%  '
Set-TextCell 'D9' '46.79'
Set-TextCell 'E9' '  +1.26%  '
Set-TextCell 'D10' '0.07932'
Set-TextCell 'E10' '  +0.86%  '
Set-TextCell 'D11' '0.9999'
Set-TextCell 'D12' '22.20'
Set-TextCell 'E12' '  +1.76%  '
Set-TextCell 'D13' '1.903.47'
Set-TextCell 'E13' '  +1.05%  '
Set-TextCell 'D14' '7.078'
Set-TextCell 'E14' '  +1.20%  '
Set-TextCell 'D15' '5.761'
Set-TextCell 'E15' '  +0.88%  '
Set-TextCell 'D16' '0.06943'
Set-TextCell 'E16' '  -0.51%  '
Set-TextCell 'D17' '88.28'
Set-TextCell 'E17' '  -0.13%  '
Set-TextCell 'E18' '  -0.04%  '
Set-TextCell 'E19' '  +0.04%  '
Set-TextCell 'E20' '  +1.90%  '
Set-TextCell 'E21' '  -0.05%  '
Set-TextCell 'D22' '29.141.15'
Set-TextCell 'E22' '  +1.27%  '
Set-TextCell 'D23' '5.358'
Set-TextCell 'E23' '  +1.47%  '
Set-TextCell 'E24' '  +0.21%  '
Set-TextCell 'D25' '2.118.11'
Set-TextCell 'E25' '  -0.44%  '
Set-TextCell 'E26' '  -2.25%  '
Set-TextCell 'D27' '156.57'
Set-TextCell 'E27' '  +2.54%  '
Set-TextCell 'D28' '19.42'
Set-TextCell 'E28' '  +0.95%  '
Set-TextCell 'D29' '6.105'
Set-TextCell 'E29' '  +4.43%  '
Set-TextCell 'D30' '1.993'
Set-TextCell 'E30' '  +0.61%  '
Set-TextCell 'D31' '118.87'
Set-TextCell 'E31' '  -0.09%  '
Set-TextCell 'D32' '0.09381'
Set-TextCell 'E32' '  +0.47%  '
Set-TextCell 'D33' '0.9243'
Set-TextCell 'E33' '  +0.17%  '
Set-TextCell 'E34' '  +0.29%  '
Set-TextCell 'E35' '  +0.66%  '
Set-TextCell 'D36' '3.271'
Set-TextCell 'E36' '  -1.58%  '
Set-TextCell 'D37' '1.200'
Set-TextCell 'E37' '  +4.07%  '
Set-TextCell 'D38' '0.05821'
Set-TextCell 'E38' '  +0.59%  '
Set-TextCell 'D39' '0.02100'
Set-TextCell 'E39' '  +1.41%  '
Set-TextCell 'D40' '7.933'
Set-TextCell 'E40' '  +3.50%  '
Set-TextCell 'E41' '  -0.13%  '
Set-TextCell 'E42' '  +1.83%  '
Set-TextCell 'E43' '  +0.72%  '
Set-TextCell 'D44' '9.911'
Set-TextCell 'E44' '  +0.69%  '
Set-TextCell 'D45' '11.92'
Set-TextCell 'E45' '  +1.14%  '
Set-TextCell 'E46' '  +2.19%  '
Set-TextCell 'D47' '2.223'
Set-TextCell 'E47' '  +4.70%  '
Set-TextCell 'D48' '0.07078'
Set-TextCell 'E48' '  -1.87%  '
Set-TextCell 'E49' '  +2.33%  '
Set-TextCell 'D50' '2.554'
Set-TextCell 'E50' '  +5.48%  '
Set-TextCell 'D51' '113.08'
Set-TextCell 'E51' '  -0.18%  '
